$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 12, 13)

# Every row in this set shares the same "Latest HO Xliff Generate Date" /
# "Latest Handoff Datetime" value (they were all regenerated together), so
# all of them move from the old timestamp to the new one.
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-25 22:21:50"
    $wsDeDe.Range("H$r").Value = "2016-08-25 22:21:50"
    $wsZhCn.Range("H$r").Value = "2016-08-25 22:21:46"
}

# Set Priority column (E) to "ht" for these handed-off records on both locale sheets.
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
